$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the current extent of the data (rows/cols already populated).
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# Insert 9 new columns in front of column B (the weekly-rank history grows to the
# left each time a new week of data is published). This shifts all existing
# columns - and any special (colored) cells inside them - nine places to the
# right, exactly like using Excel's "Insert Column" command nine times.
$ws.Columns("B:J").Insert()

# New weekly date headers for row 1, newest week first, in column order B..J.
$newDates = @("Sep_08", "Aug_25", "Aug_04", "Jul_23", "Jul_17", "Jul_07", "Jun_30", "Jun_24", "Jun_16")
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $newDates[$i]
}

# Every analyst-firm row gets a rating for each of the 9 new weeks; mark them
# all "UN" (unchanged / not yet rated) just like every other historical cell.
for ($r = 2; $r -le $lastRow; $r++) {
    for ($i = 0; $i -lt 9; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($r, $col).Value = "UN"
    }
}
